$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout) values replacing the old Strike# values in column G (rows 2-15)
$kValues = @{
    2  = 7
    3  = 6
    4  = 6
    5  = 4
    6  = 5
    7  = 5
    8  = 4
    9  = 6
    10 = 2
    11 = 5
    12 = 7
    13 = 1
    14 = 5
    15 = 6
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
